# Vervollständigung Dokumentation Arbeitsverlauf Dimitri Khodak.xlsx
# Adds a new worksheet "Fr, 06.03.2020" (a continuation of "Do, 05.03.2020")
# and fills in the closing entries of the previous day's sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Locate the last existing day sheet ("Do, 05.03.2020") ---
$lastSheet = $wb.Worksheets.Item("Do, 05.03.2020")

# --- 2. Finish up that sheet: row 23 gets the "Kostenanalyseplan" entry ---
$lastSheet.Range("C23").Value = "Optimierung des Kostenanalyseplans"
$lastSheet.Range("D23").Value = 13

# --- 3. Duplicate the sheet to create the next day's log ---
$lastSheet.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($lastSheet.Index + 1)
$newSheet.Name = "Fr, 06.03.2020"

# --- 4. Clear out the copied "today" content that doesn't belong on the new day ---

# Section 1 (Lernfeld / Dauer)
$newSheet.Range("C3").Value = 2
$newSheet.Range("D3").Value = 13

# Section 2 (Anforderungen/Inhalte)
$newSheet.Range("C5").Value = "Erstellung der Doku zu Qualitätsstandards im Bezug auf das Projekt"
$newSheet.Range("D5").Value = 13
$newSheet.Range("C6").Value = $null
$newSheet.Range("D6").Value = $null

# Section 3 (Probleme geloest / Progression)
$newSheet.Range("C10").Value = "Reflektive Einsichten und schriftliche Doku der Qualitätskriterien im Bezug auf die ISO 25010"
$newSheet.Range("D10").Value = 13

# Section 4 (Probleme aufgetreten, nicht geloest)
$newSheet.Range("C14").Value = "-"
$newSheet.Range("D14").Value = $null
$newSheet.Range("C15").Value = $null
$newSheet.Range("D15").Value = $null

# Section 5 (Hausaufgaben)
$newSheet.Range("C18").Value = "Durch hohen Zeitdruck und eine hohe Abwesenheitsrate kann das Projektziel nicht erreicht werden"
$newSheet.Range("D18").Value = "11, 12, 13"

# Section 6 (letzte Bemerkung uebernommene Zeile)
$newSheet.Range("C22").Value = "-"
$newSheet.Range("D22").Value = $null
$newSheet.Range("C23").Value = $null
$newSheet.Range("D23").Value = $null

# --- 5. Selections / active states matching the final layout ---
$lastSheet.Range("D24").Select() | Out-Null
$newSheet.Range("D6").Select() | Out-Null
$newSheet.Activate() | Out-Null
